$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New date labels for rows 201-208 (column A) ---
# Written first, in ascending order, so the new shared strings are
# appended to sharedStrings.xml in the same order as the target diff.
$ws.Cells.Item(201, 1).Value = "18 08 2020"
$ws.Cells.Item(202, 1).Value = "19 08 2020"
$ws.Cells.Item(203, 1).Value = "20 08 2020"
$ws.Cells.Item(204, 1).Value = "21 08 2020"
$ws.Cells.Item(205, 1).Value = "22 08 2020"
$ws.Cells.Item(206, 1).Value = "23 08 2020"
$ws.Cells.Item(207, 1).Value = "24 08 2020"
$ws.Cells.Item(208, 1).Value = "25 08 2020"

# --- Numeric data fills ---
# E102:E109 (previously-missing District of Columbia values)
# plus the full data rows 196-203 (columns B:D, F:BE - E stays blank)

# Row 102
$ws.Cells.Item(102, 5).Value = 0.11132914704343

# Row 103
$ws.Cells.Item(103, 5).Value = 0.10008250938702

# Row 104
$ws.Cells.Item(104, 5).Value = 0.095548433048433

# Row 105
$ws.Cells.Item(105, 5).Value = 0.11482371794872

# Row 106
$ws.Cells.Item(106, 5).Value = 0.10073339359054

# Row 107
$ws.Cells.Item(107, 5).Value = 0.09259734196709

# Row 108
$ws.Cells.Item(108, 5).Value = 0.12231240981241

# Row 109
$ws.Cells.Item(109, 5).Value = 0.11616459520871

# Row 196
$ws.Cells.Item(196, 2).Value = 0.07012987515366199
$ws.Cells.Item(196, 3).Value = 0.062135642035767
$ws.Cells.Item(196, 4).Value = 0.053072491587155
$ws.Cells.Item(196, 6).Value = 0.052899893107213
$ws.Cells.Item(196, 7).Value = 0.058911403306345
$ws.Cells.Item(196, 8).Value = 0.05669246312432
$ws.Cells.Item(196, 9).Value = 0.048601511567544
$ws.Cells.Item(196, 10).Value = 0.048834827091247
$ws.Cells.Item(196, 11).Value = 0.047230476351452
$ws.Cells.Item(196, 12).Value = 0.051771241133068
$ws.Cells.Item(196, 13).Value = 0.055204550726312
$ws.Cells.Item(196, 14).Value = 0.066075734562571
$ws.Cells.Item(196, 15).Value = 0.029596948968483
$ws.Cells.Item(196, 16).Value = 0.049790577693186
$ws.Cells.Item(196, 17).Value = 0.051546909376128
$ws.Cells.Item(196, 18).Value = 0.045827284759262
$ws.Cells.Item(196, 19).Value = 0.061332404472272
$ws.Cells.Item(196, 20).Value = 0.052420378656313
$ws.Cells.Item(196, 21).Value = 0.051386158293978
$ws.Cells.Item(196, 22).Value = 0.063903775290881
$ws.Cells.Item(196, 23).Value = 0.049934686761131
$ws.Cells.Item(196, 24).Value = 0.046696587712037
$ws.Cells.Item(196, 25).Value = 0.0497850665597
$ws.Cells.Item(196, 26).Value = 0.049345238839255
$ws.Cells.Item(196, 27).Value = 0.045016838223512
$ws.Cells.Item(196, 28).Value = 0.05395614841529
$ws.Cells.Item(196, 29).Value = 0.08177527834253701
$ws.Cells.Item(196, 30).Value = 0.06701814441437
$ws.Cells.Item(196, 31).Value = 0.059055488224207
$ws.Cells.Item(196, 32).Value = 0.055108552412082
$ws.Cells.Item(196, 33).Value = 0.049342294445628
$ws.Cells.Item(196, 34).Value = 0.06937551227852599
$ws.Cells.Item(196, 35).Value = 0.052965893135604
$ws.Cells.Item(196, 36).Value = 0.051344007157345
$ws.Cells.Item(196, 37).Value = 0.059144049650255
$ws.Cells.Item(196, 38).Value = 0.052611833093535
$ws.Cells.Item(196, 39).Value = 0.056631358562013
$ws.Cells.Item(196, 40).Value = 0.051980602518773
$ws.Cells.Item(196, 41).Value = 0.059508079292524
$ws.Cells.Item(196, 42).Value = 0.048409818263941
$ws.Cells.Item(196, 43).Value = 0.050762624903668
$ws.Cells.Item(196, 44).Value = 0.084346809998481
$ws.Cells.Item(196, 45).Value = 0.052700126918369
$ws.Cells.Item(196, 46).Value = 0.05495813538546
$ws.Cells.Item(196, 47).Value = 0.053763089485918
$ws.Cells.Item(196, 48).Value = 0.062446851919988
$ws.Cells.Item(196, 49).Value = 0.0532885157007
$ws.Cells.Item(196, 50).Value = 0.064572233348554
$ws.Cells.Item(196, 51).Value = 0.052412110204156
$ws.Cells.Item(196, 52).Value = 0.087130520932675
$ws.Cells.Item(196, 53).Value = 0.059178882942528
$ws.Cells.Item(196, 54).Value = 0.050348814368076
$ws.Cells.Item(196, 55).Value = 0.053189658278003
$ws.Cells.Item(196, 56).Value = 0.057479736608708
$ws.Cells.Item(196, 57).Value = 0.05406622656393

# Row 197
$ws.Cells.Item(197, 2).Value = 0.04996359924355
$ws.Cells.Item(197, 3).Value = 0.050973281077917
$ws.Cells.Item(197, 4).Value = 0.049439015772483
$ws.Cells.Item(197, 6).Value = 0.045918842200938
$ws.Cells.Item(197, 7).Value = 0.057675415668535
$ws.Cells.Item(197, 8).Value = 0.048006198244247
$ws.Cells.Item(197, 9).Value = 0.040191081959372
$ws.Cells.Item(197, 10).Value = 0.041828537413948
$ws.Cells.Item(197, 11).Value = 0.038506586002081
$ws.Cells.Item(197, 12).Value = 0.04188942778802
$ws.Cells.Item(197, 13).Value = 0.047948414164521
$ws.Cells.Item(197, 14).Value = 0.050637220964548
$ws.Cells.Item(197, 15).Value = 0.027209960040514
$ws.Cells.Item(197, 16).Value = 0.044892921586688
$ws.Cells.Item(197, 17).Value = 0.045491076842821
$ws.Cells.Item(197, 18).Value = 0.040737311937812
$ws.Cells.Item(197, 19).Value = 0.054065710006823
$ws.Cells.Item(197, 20).Value = 0.046721690993539
$ws.Cells.Item(197, 21).Value = 0.0465721465297
$ws.Cells.Item(197, 22).Value = 0.057865758893067
$ws.Cells.Item(197, 23).Value = 0.048677484098566
$ws.Cells.Item(197, 24).Value = 0.04495190384313
$ws.Cells.Item(197, 25).Value = 0.049409219479478
$ws.Cells.Item(197, 26).Value = 0.042081306445438
$ws.Cells.Item(197, 27).Value = 0.034783578335399
$ws.Cells.Item(197, 28).Value = 0.040875702786805
$ws.Cells.Item(197, 29).Value = 0.072821059393351
$ws.Cells.Item(197, 30).Value = 0.053794554475808
$ws.Cells.Item(197, 31).Value = 0.045007539994904
$ws.Cells.Item(197, 32).Value = 0.044863714313876
$ws.Cells.Item(197, 33).Value = 0.042121887658935
$ws.Cells.Item(197, 34).Value = 0.056896099068887
$ws.Cells.Item(197, 35).Value = 0.039209161260725
$ws.Cells.Item(197, 36).Value = 0.03708528090917
$ws.Cells.Item(197, 37).Value = 0.041609588999568
$ws.Cells.Item(197, 38).Value = 0.038915035104978
$ws.Cells.Item(197, 39).Value = 0.042329852577869
$ws.Cells.Item(197, 40).Value = 0.046492595992343
$ws.Cells.Item(197, 41).Value = 0.054827624478519
$ws.Cells.Item(197, 42).Value = 0.045825665719089
$ws.Cells.Item(197, 43).Value = 0.048112039749494
$ws.Cells.Item(197, 44).Value = 0.076990822117656
$ws.Cells.Item(197, 45).Value = 0.04856096228815
$ws.Cells.Item(197, 46).Value = 0.046220820727454
$ws.Cells.Item(197, 47).Value = 0.046451241531909
$ws.Cells.Item(197, 48).Value = 0.053589164636181
$ws.Cells.Item(197, 49).Value = 0.043410921441283
$ws.Cells.Item(197, 50).Value = 0.05122443743405
$ws.Cells.Item(197, 51).Value = 0.041271560176426
$ws.Cells.Item(197, 52).Value = 0.078598355430161
$ws.Cells.Item(197, 53).Value = 0.046739903273698
$ws.Cells.Item(197, 54).Value = 0.040661851634819
$ws.Cells.Item(197, 55).Value = 0.042342984278177
$ws.Cells.Item(197, 56).Value = 0.046752423851112
$ws.Cells.Item(197, 57).Value = 0.040729805510649

# Row 198
$ws.Cells.Item(198, 2).Value = 0.038392935214549
$ws.Cells.Item(198, 3).Value = 0.028183930485634
$ws.Cells.Item(198, 4).Value = 0.028789965391478
$ws.Cells.Item(198, 6).Value = 0.033435218174005
$ws.Cells.Item(198, 7).Value = 0.039721025564953
$ws.Cells.Item(198, 8).Value = 0.04085547227037
$ws.Cells.Item(198, 9).Value = 0.036714206903835
$ws.Cells.Item(198, 10).Value = 0.047641624251552
$ws.Cells.Item(198, 11).Value = 0.034595793720407
$ws.Cells.Item(198, 12).Value = 0.027822818322094
$ws.Cells.Item(198, 13).Value = 0.022835671553375
$ws.Cells.Item(198, 14).Value = 0.051767307352726
$ws.Cells.Item(198, 15).Value = 0.024431864968095
$ws.Cells.Item(198, 16).Value = 0.026851594872582
$ws.Cells.Item(198, 17).Value = 0.03211263516851
$ws.Cells.Item(198, 18).Value = 0.027636176680459
$ws.Cells.Item(198, 19).Value = 0.025588081721554
$ws.Cells.Item(198, 20).Value = 0.025968842418658
$ws.Cells.Item(198, 21).Value = 0.026222207525661
$ws.Cells.Item(198, 22).Value = 0.030000193895122
$ws.Cells.Item(198, 23).Value = 0.033829648636746
$ws.Cells.Item(198, 24).Value = 0.028540310748576
$ws.Cells.Item(198, 25).Value = 0.033768199103908
$ws.Cells.Item(198, 26).Value = 0.028595306732728
$ws.Cells.Item(198, 27).Value = 0.025723290945697
$ws.Cells.Item(198, 28).Value = 0.025289351079584
$ws.Cells.Item(198, 29).Value = 0.067231106763794
$ws.Cells.Item(198, 30).Value = 0.024284597483213
$ws.Cells.Item(198, 31).Value = 0.039513094530257
$ws.Cells.Item(198, 32).Value = 0.025032335538238
$ws.Cells.Item(198, 33).Value = 0.027642980584277
$ws.Cells.Item(198, 34).Value = 0.026289844306393
$ws.Cells.Item(198, 35).Value = 0.029405199245675
$ws.Cells.Item(198, 36).Value = 0.028888569597625
$ws.Cells.Item(198, 37).Value = 0.03800262249256
$ws.Cells.Item(198, 38).Value = 0.030886559791971
$ws.Cells.Item(198, 39).Value = 0.03604556017805
$ws.Cells.Item(198, 40).Value = 0.029686267481555
$ws.Cells.Item(198, 41).Value = 0.026548703536699
$ws.Cells.Item(198, 42).Value = 0.034506855888963
$ws.Cells.Item(198, 43).Value = 0.033067300996001
$ws.Cells.Item(198, 44).Value = 0.062722511512482
$ws.Cells.Item(198, 45).Value = 0.03607451775163
$ws.Cells.Item(198, 46).Value = 0.025799162529091
$ws.Cells.Item(198, 47).Value = 0.030189513384241
$ws.Cells.Item(198, 48).Value = 0.024586032797689
$ws.Cells.Item(198, 49).Value = 0.023378277418083
$ws.Cells.Item(198, 50).Value = 0.03197478712492
$ws.Cells.Item(198, 51).Value = 0.023197897611883
$ws.Cells.Item(198, 52).Value = 0.056712274811516
$ws.Cells.Item(198, 53).Value = 0.032641629434382
$ws.Cells.Item(198, 54).Value = 0.026391940252184
$ws.Cells.Item(198, 55).Value = 0.02362164418858
$ws.Cells.Item(198, 56).Value = 0.025005749223137
$ws.Cells.Item(198, 57).Value = 0.02994469247671

# Row 199
$ws.Cells.Item(199, 2).Value = 0.062513102148142
$ws.Cells.Item(199, 3).Value = 0.025205300319134
$ws.Cells.Item(199, 4).Value = 0.028073706186049
$ws.Cells.Item(199, 6).Value = 0.043406804993669
$ws.Cells.Item(199, 7).Value = 0.031807848847324
$ws.Cells.Item(199, 8).Value = 0.027180784017262
$ws.Cells.Item(199, 9).Value = 0.02465054463835
$ws.Cells.Item(199, 10).Value = 0.036157124039717
$ws.Cells.Item(199, 11).Value = 0.022206735398476
$ws.Cells.Item(199, 12).Value = 0.020702772909272
$ws.Cells.Item(199, 13).Value = 0.019906077753314
$ws.Cells.Item(199, 14).Value = 0.062138726121638
$ws.Cells.Item(199, 15).Value = 0.024483150239941
$ws.Cells.Item(199, 16).Value = 0.022949573111941
$ws.Cells.Item(199, 17).Value = 0.030024753191083
$ws.Cells.Item(199, 18).Value = 0.027560295462208
$ws.Cells.Item(199, 19).Value = 0.023646896861581
$ws.Cells.Item(199, 20).Value = 0.022525057184249
$ws.Cells.Item(199, 21).Value = 0.022590281818508
$ws.Cells.Item(199, 22).Value = 0.025828494611373
$ws.Cells.Item(199, 23).Value = 0.030267884146397
$ws.Cells.Item(199, 24).Value = 0.024740874201273
$ws.Cells.Item(199, 25).Value = 0.030166704187566
$ws.Cells.Item(199, 26).Value = 0.025997771986411
$ws.Cells.Item(199, 27).Value = 0.02350810009719
$ws.Cells.Item(199, 28).Value = 0.021672797803432
$ws.Cells.Item(199, 29).Value = 0.072236608159482
$ws.Cells.Item(199, 30).Value = 0.020165008366222
$ws.Cells.Item(199, 31).Value = 0.0318284359112
$ws.Cells.Item(199, 32).Value = 0.019334637119969
$ws.Cells.Item(199, 33).Value = 0.023160429229266
$ws.Cells.Item(199, 34).Value = 0.021336789042582
$ws.Cells.Item(199, 35).Value = 0.024176230529161
$ws.Cells.Item(199, 36).Value = 0.023708389253836
$ws.Cells.Item(199, 37).Value = 0.03164190670689
$ws.Cells.Item(199, 38).Value = 0.026095950641195
$ws.Cells.Item(199, 39).Value = 0.029967681319525
$ws.Cells.Item(199, 40).Value = 0.023259198453475
$ws.Cells.Item(199, 41).Value = 0.019468982862262
$ws.Cells.Item(199, 42).Value = 0.025199467165631
$ws.Cells.Item(199, 43).Value = 0.025603031963977
$ws.Cells.Item(199, 44).Value = 0.071356641519062
$ws.Cells.Item(199, 45).Value = 0.030412652731413
$ws.Cells.Item(199, 46).Value = 0.021322983626246
$ws.Cells.Item(199, 47).Value = 0.026277948859102
$ws.Cells.Item(199, 48).Value = 0.021070774615129
$ws.Cells.Item(199, 49).Value = 0.02306419703959
$ws.Cells.Item(199, 50).Value = 0.032150557298737
$ws.Cells.Item(199, 51).Value = 0.027228075017764
$ws.Cells.Item(199, 52).Value = 0.062277584608193
$ws.Cells.Item(199, 53).Value = 0.040038507580829
$ws.Cells.Item(199, 54).Value = 0.031423155386562
$ws.Cells.Item(199, 55).Value = 0.027569347101493
$ws.Cells.Item(199, 56).Value = 0.030757441515313
$ws.Cells.Item(199, 57).Value = 0.033355564714437

# Row 200
$ws.Cells.Item(200, 2).Value = 0.072960352701781
$ws.Cells.Item(200, 3).Value = 0.06932199811913101
$ws.Cells.Item(200, 4).Value = 0.066925731558018
$ws.Cells.Item(200, 6).Value = 0.064785732025993
$ws.Cells.Item(200, 7).Value = 0.05400737072254
$ws.Cells.Item(200, 8).Value = 0.058101476716908
$ws.Cells.Item(200, 9).Value = 0.047033650437939
$ws.Cells.Item(200, 10).Value = 0.045393671958731
$ws.Cells.Item(200, 11).Value = 0.046543056329954
$ws.Cells.Item(200, 12).Value = 0.048302191319422
$ws.Cells.Item(200, 13).Value = 0.056917632671693
$ws.Cells.Item(200, 14).Value = 0.059893740003529
$ws.Cells.Item(200, 15).Value = 0.031863173165427
$ws.Cells.Item(200, 16).Value = 0.057071495050643
$ws.Cells.Item(200, 17).Value = 0.052534659870061
$ws.Cells.Item(200, 18).Value = 0.050837481511949
$ws.Cells.Item(200, 19).Value = 0.066616509215547
$ws.Cells.Item(200, 20).Value = 0.05611015856948
$ws.Cells.Item(200, 21).Value = 0.055679341908934
$ws.Cells.Item(200, 22).Value = 0.065804699254736
$ws.Cells.Item(200, 23).Value = 0.046187757553233
$ws.Cells.Item(200, 24).Value = 0.046487281825346
$ws.Cells.Item(200, 25).Value = 0.05052298232608
$ws.Cells.Item(200, 26).Value = 0.04566867346136
$ws.Cells.Item(200, 27).Value = 0.042251847058997
$ws.Cells.Item(200, 28).Value = 0.051995135020124
$ws.Cells.Item(200, 29).Value = 0.074466130726228
$ws.Cells.Item(200, 30).Value = 0.066568580994567
$ws.Cells.Item(200, 31).Value = 0.055043136037183
$ws.Cells.Item(200, 32).Value = 0.057918073458903
$ws.Cells.Item(200, 33).Value = 0.049692493388909
$ws.Cells.Item(200, 34).Value = 0.07141605962481
$ws.Cells.Item(200, 35).Value = 0.048335461134739
$ws.Cells.Item(200, 36).Value = 0.046267980589534
$ws.Cells.Item(200, 37).Value = 0.052352553003399
$ws.Cells.Item(200, 38).Value = 0.048007423180615
$ws.Cells.Item(200, 39).Value = 0.051069539854303
$ws.Cells.Item(200, 40).Value = 0.05190885976483
$ws.Cells.Item(200, 41).Value = 0.05779362027516
$ws.Cells.Item(200, 42).Value = 0.044748121724452
$ws.Cells.Item(200, 43).Value = 0.0471100035955
$ws.Cells.Item(200, 44).Value = 0.08128173235582099
$ws.Cells.Item(200, 45).Value = 0.047820518562304
$ws.Cells.Item(200, 46).Value = 0.053774842675955
$ws.Cells.Item(200, 47).Value = 0.049901073420962
$ws.Cells.Item(200, 48).Value = 0.058455658756022
$ws.Cells.Item(200, 49).Value = 0.049161755806869
$ws.Cells.Item(200, 50).Value = 0.060912394501243
$ws.Cells.Item(200, 51).Value = 0.051714499308891
$ws.Cells.Item(200, 52).Value = 0.085531671329859
$ws.Cells.Item(200, 53).Value = 0.051947562925779
$ws.Cells.Item(200, 54).Value = 0.046333157586519
$ws.Cells.Item(200, 55).Value = 0.052192904900678
$ws.Cells.Item(200, 56).Value = 0.055791360918829
$ws.Cells.Item(200, 57).Value = 0.054232892164952

# Row 201
$ws.Cells.Item(201, 2).Value = 0.062942167715434
$ws.Cells.Item(201, 3).Value = 0.076832107318223
$ws.Cells.Item(201, 4).Value = 0.055678755166565
$ws.Cells.Item(201, 6).Value = 0.047655865065464
$ws.Cells.Item(201, 7).Value = 0.055678611048956
$ws.Cells.Item(201, 8).Value = 0.062417078703609
$ws.Cells.Item(201, 9).Value = 0.04988538520012
$ws.Cells.Item(201, 10).Value = 0.047894608505466
$ws.Cells.Item(201, 11).Value = 0.049321925610216
$ws.Cells.Item(201, 12).Value = 0.054629376879733
$ws.Cells.Item(201, 13).Value = 0.061014111405554
$ws.Cells.Item(201, 14).Value = 0.057770902372154
$ws.Cells.Item(201, 15).Value = 0.034628384945807
$ws.Cells.Item(201, 16).Value = 0.056159319742915
$ws.Cells.Item(201, 17).Value = 0.056507359729599
$ws.Cells.Item(201, 18).Value = 0.051022536519295
$ws.Cells.Item(201, 19).Value = 0.06453060423691601
$ws.Cells.Item(201, 20).Value = 0.057236967645642
$ws.Cells.Item(201, 21).Value = 0.058539592324018
$ws.Cells.Item(201, 22).Value = 0.07008056554825
$ws.Cells.Item(201, 23).Value = 0.056939755506581
$ws.Cells.Item(201, 24).Value = 0.052582901909005
$ws.Cells.Item(201, 25).Value = 0.056920063971646
$ws.Cells.Item(201, 26).Value = 0.05748374272499
$ws.Cells.Item(201, 27).Value = 0.052937430109873
$ws.Cells.Item(201, 28).Value = 0.060253640813648
$ws.Cells.Item(201, 29).Value = 0.06532111667827099
$ws.Cells.Item(201, 30).Value = 0.076921930576588
$ws.Cells.Item(201, 31).Value = 0.06502790928730701
$ws.Cells.Item(201, 32).Value = 0.059095138491728
$ws.Cells.Item(201, 33).Value = 0.056943377445941
$ws.Cells.Item(201, 34).Value = 0.081050678382476
$ws.Cells.Item(201, 35).Value = 0.052352867500923
$ws.Cells.Item(201, 36).Value = 0.050634413198851
$ws.Cells.Item(201, 37).Value = 0.057325487631923
$ws.Cells.Item(201, 38).Value = 0.054893076233475
$ws.Cells.Item(201, 39).Value = 0.052194866831231
$ws.Cells.Item(201, 40).Value = 0.05529994778309
$ws.Cells.Item(201, 41).Value = 0.063391427171613
$ws.Cells.Item(201, 42).Value = 0.04988516557589
$ws.Cells.Item(201, 43).Value = 0.052007948535185
$ws.Cells.Item(201, 44).Value = 0.07800664054602199
$ws.Cells.Item(201, 45).Value = 0.05185515844966
$ws.Cells.Item(201, 46).Value = 0.057796083385071
$ws.Cells.Item(201, 47).Value = 0.057263617322551
$ws.Cells.Item(201, 48).Value = 0.063257931243629
$ws.Cells.Item(201, 49).Value = 0.052064438965304
$ws.Cells.Item(201, 50).Value = 0.060404657451295
$ws.Cells.Item(201, 51).Value = 0.048999218105517
$ws.Cells.Item(201, 52).Value = 0.075580539854932
$ws.Cells.Item(201, 53).Value = 0.055060141615547
$ws.Cells.Item(201, 54).Value = 0.045654863624054
$ws.Cells.Item(201, 55).Value = 0.050292260535609
$ws.Cells.Item(201, 56).Value = 0.052228961360273
$ws.Cells.Item(201, 57).Value = 0.053223387671244

# Row 202
$ws.Cells.Item(202, 2).Value = 0.057239871530782
$ws.Cells.Item(202, 3).Value = 0.06447276074363301
$ws.Cells.Item(202, 4).Value = 0.05357039171125
$ws.Cells.Item(202, 6).Value = 0.047652040544
$ws.Cells.Item(202, 7).Value = 0.050575683997647
$ws.Cells.Item(202, 8).Value = 0.059776391261376
$ws.Cells.Item(202, 9).Value = 0.046825413423194
$ws.Cells.Item(202, 10).Value = 0.042901333850918
$ws.Cells.Item(202, 11).Value = 0.046987497439987
$ws.Cells.Item(202, 12).Value = 0.049305388771167
$ws.Cells.Item(202, 13).Value = 0.056237667852631
$ws.Cells.Item(202, 14).Value = 0.056333319587841
$ws.Cells.Item(202, 15).Value = 0.028768238669416
$ws.Cells.Item(202, 16).Value = 0.049721410423936
$ws.Cells.Item(202, 17).Value = 0.055050115481605
$ws.Cells.Item(202, 18).Value = 0.047251069500908
$ws.Cells.Item(202, 19).Value = 0.058527510373764
$ws.Cells.Item(202, 20).Value = 0.049550531892891
$ws.Cells.Item(202, 21).Value = 0.04821751487933
$ws.Cells.Item(202, 22).Value = 0.057925270552898
$ws.Cells.Item(202, 23).Value = 0.039910539896107
$ws.Cells.Item(202, 24).Value = 0.03911811032908
$ws.Cells.Item(202, 25).Value = 0.043547176624838
$ws.Cells.Item(202, 26).Value = 0.04087832212925
$ws.Cells.Item(202, 27).Value = 0.038187359763098
$ws.Cells.Item(202, 28).Value = 0.04965268283228
$ws.Cells.Item(202, 29).Value = 0.07377829750591799
$ws.Cells.Item(202, 30).Value = 0.063011155551582
$ws.Cells.Item(202, 31).Value = 0.052802371570466
$ws.Cells.Item(202, 32).Value = 0.054908393115583
$ws.Cells.Item(202, 33).Value = 0.053549670043582
$ws.Cells.Item(202, 34).Value = 0.070975850127259
$ws.Cells.Item(202, 35).Value = 0.045216950484251
$ws.Cells.Item(202, 36).Value = 0.042959740766357
$ws.Cells.Item(202, 37).Value = 0.050057455009089
$ws.Cells.Item(202, 38).Value = 0.043987104824909
$ws.Cells.Item(202, 39).Value = 0.047445989511853
$ws.Cells.Item(202, 40).Value = 0.052334976753218
$ws.Cells.Item(202, 41).Value = 0.062217357868657
$ws.Cells.Item(202, 42).Value = 0.048551612478626
$ws.Cells.Item(202, 43).Value = 0.051241657679244
$ws.Cells.Item(202, 44).Value = 0.08112126482241799
$ws.Cells.Item(202, 45).Value = 0.050626496979557
$ws.Cells.Item(202, 46).Value = 0.058890726616584
$ws.Cells.Item(202, 47).Value = 0.06665741362137501
$ws.Cells.Item(202, 48).Value = 0.06372723327755
$ws.Cells.Item(202, 49).Value = 0.054972815406301
$ws.Cells.Item(202, 50).Value = 0.06539878111279999
$ws.Cells.Item(202, 51).Value = 0.050645384743208
$ws.Cells.Item(202, 52).Value = 0.07480124795177399
$ws.Cells.Item(202, 53).Value = 0.058554655005836
$ws.Cells.Item(202, 54).Value = 0.047749634969037
$ws.Cells.Item(202, 55).Value = 0.052704234724526
$ws.Cells.Item(202, 56).Value = 0.055726604373504
$ws.Cells.Item(202, 57).Value = 0.057473523415556

# Row 203
$ws.Cells.Item(203, 2).Value = 0.053787858266768
$ws.Cells.Item(203, 3).Value = 0.064217578083422
$ws.Cells.Item(203, 4).Value = 0.053750320542592
$ws.Cells.Item(203, 6).Value = 0.04605977530563
$ws.Cells.Item(203, 7).Value = 0.04841089934952
$ws.Cells.Item(203, 8).Value = 0.054910051010528
$ws.Cells.Item(203, 9).Value = 0.045298013551403
$ws.Cells.Item(203, 10).Value = 0.045568723684285
$ws.Cells.Item(203, 11).Value = 0.045344160726492
$ws.Cells.Item(203, 12).Value = 0.049803617021394
$ws.Cells.Item(203, 13).Value = 0.057224339272908
$ws.Cells.Item(203, 14).Value = 0.058201668711171
$ws.Cells.Item(203, 15).Value = 0.03096431500829
$ws.Cells.Item(203, 16).Value = 0.058321610316679
$ws.Cells.Item(203, 17).Value = 0.057980425230833
$ws.Cells.Item(203, 18).Value = 0.051512846672651
$ws.Cells.Item(203, 19).Value = 0.06444690943535
$ws.Cells.Item(203, 20).Value = 0.058521902570185
$ws.Cells.Item(203, 21).Value = 0.056540555250892
$ws.Cells.Item(203, 22).Value = 0.06893271803310901
$ws.Cells.Item(203, 23).Value = 0.050147023546503
$ws.Cells.Item(203, 24).Value = 0.048741605684827
$ws.Cells.Item(203, 25).Value = 0.053376079052184
$ws.Cells.Item(203, 26).Value = 0.047759445023717
$ws.Cells.Item(203, 27).Value = 0.042431966817826
$ws.Cells.Item(203, 28).Value = 0.053749584724794
$ws.Cells.Item(203, 29).Value = 0.07613043534031499
$ws.Cells.Item(203, 30).Value = 0.068119494243815
$ws.Cells.Item(203, 31).Value = 0.054885893351456
$ws.Cells.Item(203, 32).Value = 0.059547309180769
$ws.Cells.Item(203, 33).Value = 0.060599782360616
$ws.Cells.Item(203, 34).Value = 0.076413115947528
$ws.Cells.Item(203, 35).Value = 0.047284140553018
$ws.Cells.Item(203, 36).Value = 0.043551571663421
$ws.Cells.Item(203, 37).Value = 0.049339181288424
$ws.Cells.Item(203, 38).Value = 0.045164912472906
$ws.Cells.Item(203, 39).Value = 0.048563201349032
$ws.Cells.Item(203, 40).Value = 0.05880125354473
$ws.Cells.Item(203, 41).Value = 0.06868179999799801
$ws.Cells.Item(203, 42).Value = 0.053728862696027
$ws.Cells.Item(203, 43).Value = 0.055836698956009
$ws.Cells.Item(203, 44).Value = 0.06586051743850201
$ws.Cells.Item(203, 45).Value = 0.054252526864763
$ws.Cells.Item(203, 46).Value = 0.060289487935493
$ws.Cells.Item(203, 47).Value = 0.072669125740899
$ws.Cells.Item(203, 48).Value = 0.06542831531186399
$ws.Cells.Item(203, 49).Value = 0.055251747017938
$ws.Cells.Item(203, 50).Value = 0.061065086354791
$ws.Cells.Item(203, 51).Value = 0.049918297233704
$ws.Cells.Item(203, 52).Value = 0.057547184242071
$ws.Cells.Item(203, 53).Value = 0.052621945035247
$ws.Cells.Item(203, 54).Value = 0.045407570113341
$ws.Cells.Item(203, 55).Value = 0.049721665435199
$ws.Cells.Item(203, 56).Value = 0.054491847067095
$ws.Cells.Item(203, 57).Value = 0.054320330233501
